$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$updates1 = @(
    ,(3, 5063)
    ,(5, 7336)
    ,(11, 23)
    ,(12, 4279)
    ,(13, 1731)
    ,(15, 97)
    ,(16, 2882)
    ,(18, 562)
    ,(20, 472)
    ,(21, 424)
    ,(22, 448)
    ,(23, 284)
    ,(24, 88)
    ,(25, 1673)
    ,(26, 1158)
    ,(28, 1358)
    ,(30, 574)
    ,(31, 21)
    ,(33, 23)
    ,(34, 56)
    ,(35, 103)
    ,(36, 51)
    ,(37, 2788)
    ,(38, 695)
    ,(39, 18)
    ,(40, 39)
)
foreach ($pair in $updates1) {
    $ws.Cells.Item($pair[0], 6).Value = $pair[1]
}

$ws = $wb.Worksheets.Item(2)
$updates2 = @(
    ,(2, 5)
)
foreach ($pair in $updates2) {
    $ws.Cells.Item($pair[0], 6).Value = $pair[1]
}

$ws = $wb.Worksheets.Item(4)
$updates4 = @(
    ,(3, 5063)
    ,(5, 7336)
    ,(11, 23)
    ,(12, 4279)
    ,(13, 1731)
    ,(15, 97)
    ,(16, 2882)
    ,(18, 562)
    ,(20, 472)
    ,(21, 424)
    ,(22, 448)
    ,(23, 284)
    ,(24, 88)
    ,(25, 1673)
    ,(26, 1158)
    ,(28, 1358)
    ,(30, 574)
    ,(31, 21)
    ,(33, 23)
    ,(34, 56)
    ,(35, 103)
    ,(36, 51)
    ,(37, 2788)
    ,(38, 5)
    ,(39, 695)
    ,(40, 18)
    ,(41, 39)
)
foreach ($pair in $updates4) {
    $ws.Cells.Item($pair[0], 6).Value = $pair[1]
}
